$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row cells: "_old" -> "_FV2404", "_new" -> "_FV2410"
$headers = @(
    "Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404",
    "Segment ID_FV2404","Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404","Bedingung_FV2404","diff",
    "Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410",
    "Segment ID_FV2410","Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410","Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2. Convert the existing range into an Excel Table (ListObject) named "Table1"
$range = $ws.Range("A1:U89")
$list = $ws.ListObjects.Add(1, $range, $null, 1)
$list.Name = "Table1"
$list.TableStyle = ""

# 3. Freeze the header row (split after row 1)
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
